$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Drop the stale _xlchart.v1.* hidden defined names (leftovers from
#    a previous chart-insertion workflow that no longer apply).
# ------------------------------------------------------------------
$namesToDelete = @()
foreach ($n in $wb.Names) {
    if ($n.Name -like "_xlchart.v1.*") {
        $namesToDelete += $n.Name
    }
}
foreach ($nm in $namesToDelete) {
    $wb.Names.Item($nm).Delete()
}

# ------------------------------------------------------------------
# 2. Move the "Data" sheet's selection (cosmetic, matches authored file).
# ------------------------------------------------------------------
$dataSheet = $wb.Worksheets.Item("Data")
$dataSheet.Range("H24").Select()

# ------------------------------------------------------------------
# 3. Add the new "CWE" worksheet as the last tab in the workbook.
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "CWE"

# ------------------------------------------------------------------
# 4. Write the year header row (row 1) and pre-format the whole A1:N1
#    band to be center-aligned, then merge each year's two columns.
# ------------------------------------------------------------------
$years = @(2011, 2012, 2013, 2014, 2015, 2016, 2017)
$yearCols = @(1, 3, 5, 7, 9, 11, 13)

$ws.Range("A1:N1").HorizontalAlignment = -4108

for ($i = 0; $i -lt $years.Length; $i++) {
    $ws.Cells.Item(1, $yearCols[$i]).Value = $years[$i]
}

$mergeRanges = @("A1:B1", "C1:D1", "E1:F1", "G1:H1", "I1:J1", "K1:L1", "M1:N1")
foreach ($mr in $mergeRanges) {
    $ws.Range($mr).Merge()
}

# ------------------------------------------------------------------
# 5. Write the CWE counts (numeric columns) and CWE names (label
#    columns) for each of the 7 years, 10 rows of data each.
# ------------------------------------------------------------------
$countCols = @(1, 3, 5, 7, 9, 11, 13)
$labelCols = @(2, 4, 6, 8, 10, 12, 14)

$counts = @(
    @(728, 536, 478, 451, 398, 365, 321, 208, 162, 149),
    @(844, 784, 680, 594, 374, 257, 236, 229, 221, 155),
    @(960, 774, 697, 634, 525, 292, 261, 169, 168, 160),
    @(1547, 948, 822, 709, 689, 509, 418, 349, 270, 269),
    @(1116, 776, 701, 648, 647, 577, 545, 352, 230, 221),
    @(1363, 846, 788, 744, 689, 652, 584, 304, 217, 205),
    @(2322, 1274, 1208, 1175, 1023, 967, 459, 374, 355, 321)
)

$labels = @(
    @("CWE-119", "NVD-CWE-noinfo", "CWE-79", "CWE-20", "CWE-399", "CWE-264", "CWE-200", "NVD-CWE-Other", "CWE-89", "CWE-189"),
    @("NVD-CWE-noinfo", "CWE-79", "CWE-119", "CWE-264", "CWE-20", "CWE-399", "NVD-CWE-Other", "CWE-200", "CWE-89", "CWE-352"),
    @("NVD-CWE-noinfo", "CWE-119", "CWE-79", "CWE-264", "CWE-20", "CWE-399", "CWE-200", "CWE-89", "CWE-352", "CWE-94"),
    @("CWE-310", "CWE-79", "CWE-119", "NVD-CWE-noinfo", "CWE-264", "CWE-20", "CWE-200", "NVD-CWE-Other", "CWE-89", "CWE-399"),
    @("CWE-119", "CWE-79", "CWE-200", "CWE-264", "NVD-CWE-noinfo", "NVD-CWE-Other", "CWE-20", "CWE-399", "CWE-352", "CWE-89"),
    @("CWE-119", "CWE-200", "NVD-CWE-noinfo", "CWE-264", "CWE-79", "CWE-20", "CWE-284", "NVD-CWE-Other", "CWE-399", "CWE-310"),
    @("CWE-119", "CWE-79", "CWE-200", "CWE-284", "CWE-264", "CWE-20", "CWE-89", "CWE-125", "CWE-399", "CWE-22")
)

for ($col = 0; $col -lt 7; $col++) {
    for ($row = 0; $row -lt 10; $row++) {
        $sheetRow = $row + 2
        $ws.Cells.Item($sheetRow, $countCols[$col]).Value = $counts[$col][$row]
        $ws.Cells.Item($sheetRow, $labelCols[$col]).Value = $labels[$col][$row]
    }
}

# ------------------------------------------------------------------
# 6. Cosmetic sheet-level bits: portrait print orientation and final
#    selection so "CWE" ends up the active/selected tab.
# ------------------------------------------------------------------
$ws.PageSetup.Orientation = 1
$ws.Range("K15").Select()
